$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plot")

# Insert a new row at position 22, pushing existing rows (22-33) down to (23-34).
# This also shifts the frozen pane selection & keeps per-row styles consistent
# with the row that previously occupied position 22 (since Insert copies format
# from the row above, which already carries the same style pattern as old row22).
$ws.Rows.Item(22).Insert()

# Fill in the data for the newly inserted row 22 (Wavelets V2 Coiflet2 4 Level
# T2s "with Hampel Filter" result row).
$ws.Range("A22").Value = "WaveletsV2 (Coiflet2,  4 Level, T 2s) with Hampel Filter"
$ws.Range("B22").Value = 0.67100000000000004
$ws.Range("C22").Value = 0.89700000000000002
$ws.Range("D22").Value = 0.875
$ws.Range("E22").Value = 0.78100000000000003
$ws.Range("F22").Value = 0.89600000000000002
$ws.Range("G22").Value = 0.92900000000000005
$ws.Range("H22").Value = 0.79
$ws.Range("I22").Value = 0.84399999999999997
$ws.Range("J22").Value = 0.52
$ws.Range("K22").Value = 0.85899999999999999
$ws.Range("L22").Value = 0.81100000000000005
$ws.Range("M22").Value = 0.76500000000000001
$ws.Range("N22").Value = 0.81299999999999994
$ws.Range("Q22").Value = 0.93899999999999995

# Extend the conditional formatting ranges by one row, to account for the
# newly inserted row (they previously covered up to row 24/25, now 25/26).
# Re-point the existing rules at their new ranges (rather than deleting and
# recreating them) so the dxfId/colors/rule grouping are preserved exactly.
$fcsOld24 = $ws.Range("B2:Q24").FormatConditions
$fcsOld24.Item(1).ModifyAppliesToRange($ws.Range("B2:Q25"))

$fcsOld25 = $ws.Range("B2:Q25").FormatConditions
$fcsOld25.Item($fcsOld25.Count).ModifyAppliesToRange($ws.Range("B2:Q26"))

# Restore the selection (now shifted one row/column further due to the insert).
$ws.Range("R28").Select()

# Window geometry cosmetic update (position/size of the Excel window).
$w = $excel.ActiveWindow
$w.Left = 2440
$w.Top = 720
$w.Width = 21720
$w.Height = 13380
